$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to short machine-friendly names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the Spanish connector words (de/del/las/los/el/la/y) in
# state and municipality names, e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga"
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B10").Value = "Playas De Rosarito"
$ws.Range("B29").Value = "Benemérito De Las Américas"
$ws.Range("B52").Value = "Salto De Agua"
$ws.Range("B53").Value = "San Cristóbal De Las Casas"
$ws.Range("B77").Value = "Hidalgo Del Parral"
$ws.Range("B97").Value = "San Juan De Sabinas"
$ws.Range("A106").Value = "Ciudad De México"
$ws.Range("B129").Value = "Pánuco De Coronado"
$ws.Range("B132").Value = "San Juan Del Río"
$ws.Range("B133").Value = "San Pedro Del Gallo"
$ws.Range("A137").Value = "Estado De México"
$ws.Range("B137").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B140").Value = "Almoloya De Alquisiras"
$ws.Range("B142").Value = "Atizapán De Zaragoza"
$ws.Range("B144").Value = "Coacalco De Berriozábal"
$ws.Range("B146").Value = "Ecatepec De Morelos"
$ws.Range("B149").Value = "Naucalpan De Juárez"
$ws.Range("B152").Value = "San Felipe Del Progreso"
$ws.Range("B158").Value = "Tenango Del Valle"
$ws.Range("B161").Value = "Tlalnepantla De Baz"
$ws.Range("B164").Value = "Valle De Bravo"
$ws.Range("B165").Value = "Villa De Allende"
$ws.Range("B170").Value = "San Miguel De Allende"
$ws.Range("B171").Value = "Apaseo El Alto"
$ws.Range("B172").Value = "Apaseo El Grande"
$ws.Range("B178").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B188").Value = "San Diego De La Unión"
$ws.Range("B191").Value = "San Luis De La Paz"
$ws.Range("B192").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B193").Value = "Silao De La Victoria"
$ws.Range("B195").Value = "Valle De Santiago"
$ws.Range("B201").Value = "Acapulco De Juárez"
$ws.Range("B203").Value = "Ajuchitlán Del Progreso"
$ws.Range("B206").Value = "Atenango Del Río"
$ws.Range("B208").Value = "Atoyac De Álvarez"
$ws.Range("B209").Value = "Ayutla De Los Libres"
$ws.Range("B211").Value = "Chilapa De Álvarez"
$ws.Range("B212").Value = "Chilpancingo De Los Bravo"
$ws.Range("B217").Value = "Coyuca De Benítez"
$ws.Range("B218").Value = "Coyuca De Catalán"
$ws.Range("B220").Value = "Cuetzala Del Progreso"
$ws.Range("B221").Value = "Cutzamala De Pinzón"
$ws.Range("B225").Value = "Huitzuco De Los Figueroa"
$ws.Range("B226").Value = "Iguala De La Independencia"
$ws.Range("B227").Value = "Zihuatanejo De Azueta"
$ws.Range("B235").Value = "Taxco De Alarcón"
$ws.Range("B237").Value = "Técpan De Galeana"
$ws.Range("B240").Value = "Tixtla De Guerrero"
$ws.Range("B242").Value = "Tlapa De Comonfort"
$ws.Range("B251").Value = "Atotonilco El Grande"
$ws.Range("B256").Value = "Cuautepec De Hinojosa"
$ws.Range("B259").Value = "Huejutla De Reyes"
$ws.Range("B262").Value = "Jacala De Ledezma"
$ws.Range("B268").Value = "Nopala De Villagrán"
$ws.Range("B269").Value = "Pachuca De Soto"
$ws.Range("B272").Value = "Progreso De Obregón"
$ws.Range("B274").Value = "Santiago De Anaya"
$ws.Range("B275").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B278").Value = "Tenango De Doria"
$ws.Range("B280").Value = "Tepehuacán De Guerrero"
$ws.Range("B281").Value = "Tezontepec De Aldama"
$ws.Range("B285").Value = "Tula De Allende"
$ws.Range("B286").Value = "Tulancingo De Bravo"
$ws.Range("B290").Value = "Atotonilco El Alto"
$ws.Range("B296").Value = "Encarnación De Díaz"
$ws.Range("B298").Value = "Ixtlahuacán Del Río"
$ws.Range("B301").Value = "Lagos De Moreno"
$ws.Range("B305").Value = "Ojuelos De Jalisco"
$ws.Range("B308").Value = "Tamazula De Gordiano"
$ws.Range("B310").Value = "Tepatitlán De Morelos"
$ws.Range("B311").Value = "Tizapán El Alto"
$ws.Range("B315").Value = "Yahualica De González Gallo"
$ws.Range("B353").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B368").Value = "Coatlán Del Río"
$ws.Range("B380").Value = "Zacualpan De Amilpas"
$ws.Range("B387").Value = "Santa María Del Oro"
$ws.Range("B402").Value = "San Nicolás De Los Garza"
$ws.Range("B406").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B411").Value = "Cuilápam De Guerrero"
$ws.Range("B412").Value = "Guevea De Humboldt"
$ws.Range("B413").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B414").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B416").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B417").Value = "Oaxaca De Juárez"
$ws.Range("B446").Value = "Santa Inés Del Monte"
$ws.Range("B449").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B453").Value = "Santo Domingo De Morelos"
$ws.Range("B458").Value = "Tataltepec De Valdés"
$ws.Range("B459").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B460").Value = "Totontepec Villa De Morelos"
$ws.Range("B461").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B462").Value = "Villa De Zaachila"
$ws.Range("B477").Value = "Huehuetlán El Chico"
$ws.Range("B481").Value = "Izúcar De Matamoros"
$ws.Range("B486").Value = "Palmar De Bravo"
$ws.Range("B491").Value = "San Salvador El Seco"
$ws.Range("B492").Value = "Tecali De Herrera"
$ws.Range("B495").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B497").Value = "Tepexi De Rodríguez"
$ws.Range("B498").Value = "Tetela De Ocampo"
$ws.Range("B500").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B515").Value = "Amealco De Bonfil"
$ws.Range("B517").Value = "Cadereyta De Montes"
$ws.Range("B522").Value = "Jalpan De Serra"
$ws.Range("B523").Value = "Landa De Matamoros"
$ws.Range("B525").Value = "Pinal De Amoles"
$ws.Range("B527").Value = "San Juan Del Río"
$ws.Range("B537").Value = "Ciudad Del Maíz"
$ws.Range("B546").Value = "Mexquitic De Carmona"
$ws.Range("B552").Value = "San Ciro De Acosta"
$ws.Range("B558").Value = "Santa María Del Río"
$ws.Range("B563").Value = "Tanquián De Escobedo"
$ws.Range("B567").Value = "Villa De Arista"
$ws.Range("B568").Value = "Villa De Arriaga"
$ws.Range("B569").Value = "Villa De Guadalupe"
$ws.Range("B570").Value = "Villa De La Paz"
$ws.Range("B571").Value = "Villa De Ramos"
$ws.Range("B572").Value = "Villa De Reyes"
$ws.Range("B600").Value = "Jalpa De Méndez"
$ws.Range("B625").Value = "Soto La Marina"
$ws.Range("B646").Value = "Amatlán De Los Reyes"
$ws.Range("B651").Value = "Castillo De Teayo"
$ws.Range("B665").Value = "Hueyapan De Ocampo"
$ws.Range("B666").Value = "Ignacio De La Llave"
$ws.Range("B678").Value = "Martínez De La Torre"
$ws.Range("B684").Value = "Ozuluama De Mascareñas"
$ws.Range("B688").Value = "Paso Del Macho"
$ws.Range("B690").Value = "Poza Rica De Hidalgo"
$ws.Range("B695").Value = "Sayula De Alemán"
$ws.Range("B696").Value = "Soledad De Doblado"
$ws.Range("B713").Value = "Vega De Alatorre"
$ws.Range("B725").Value = "Concepción Del Oro"
$ws.Range("B736").Value = "Nochistlán De Mejía"
$ws.Range("B746").Value = "Villa De Cos"

# Drop the trailing footnote/source rows (754-758); the table now ends at row 752
$ws.Range("A754:D758").ClearContents()

